# 200V_Geschw.xlsx - "Fast fertig, kleine Unklarheit im zweiten Teil der Auswertung"
#
# - Turns the per-row "=0.5/B.." / "=0.5/C.." formulas in columns E and F
#   (rows 2:31) into a single filled-down formula (Excel stores these as
#   shared formulas once they are entered as one fill operation).
# - Adds a new column L with "=K./I." for every group-header row
#   (2,5,8,11,14,17,20,23,26,29) - the new "2v0" ratio the author was
#   unsure about ("kleine Unklarheit").
# - Updates dimension / selection to reflect the newly used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E/F: re-enter as one fill so Excel groups them as shared formulas ---
$ws.Range("E2:E31").Formula = "=0.5/B2"
$ws.Range("F2:F31").Formula = "=0.5/C2"

# --- new column L: ratio 2*v0_fall / (v0_auf - v0_ab) for each data group ---
# (column L has no explicit <col> style, so drop the inherited number format
# to keep these cells on the workbook's default "General" style, like Excel
# does for a brand-new column typed into next to a styled one)
$groupRows = @(2, 5, 8, 11, 14, 17, 20, 23, 26, 29)
foreach ($r in $groupRows) {
    $ws.Cells.Item($r, 12).Formula = "=K$r/I$r"
    $ws.Cells.Item($r, 12).ClearFormats()
}

# --- cosmetic follow-ups that Excel performs when the used range grows ---
$ws.Range("L30").Select()
